$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 5889
$ws.Range("F5").Value = 5889
$ws.Range("F7").Value = 2937
$ws.Range("F9").Value = 397
$ws.Range("F10").Value = 433
$ws.Range("F11").Value = 108
$ws.Range("F13").Value = 688
$ws.Range("F14").Value = 186
$ws.Range("F15").Value = 4235
$ws.Range("F16").Value = 4235
$ws.Range("F17").Value = 95
$ws.Range("F19").Value = 100
$ws.Range("F21").Value = 192
$ws.Range("F22").Value = 59
$ws.Range("F23").Value = 6370
$ws.Range("F24").Value = 6370
$ws.Range("F25").Value = 223
$ws.Range("F26").Value = 92
$ws.Range("F31").Value = 5761
$ws.Range("F32").Value = 1619
$ws.Range("F35").Value = 5900
$ws.Range("F36").Value = 101
$ws.Range("F39").Value = 77
$ws.Range("F40").Value = 268
$ws.Range("F41").Value = 4002
$ws.Range("F42").Value = 182
$ws.Range("F44").Value = 15
$ws.Range("F45").Value = 2394
$ws.Range("F48").Value = 1005
$ws.Range("F49").Value = 14
$ws.Range("F50").Value = 289
$ws.Range("F51").Value = 2033
$ws.Range("F52").Value = 16

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 125
$ws.Range("F4").Value = 25
$ws.Range("F5").Value = 91

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1406

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1406
$ws.Range("F4").Value = 5889
$ws.Range("F5").Value = 5889
$ws.Range("F7").Value = 2937
$ws.Range("F9").Value = 433
$ws.Range("F10").Value = 108
$ws.Range("F13").Value = 186
$ws.Range("F14").Value = 4235
$ws.Range("F15").Value = 4235
$ws.Range("F16").Value = 95
$ws.Range("F18").Value = 100
$ws.Range("F20").Value = 192
$ws.Range("F21").Value = 59
$ws.Range("F22").Value = 6370
$ws.Range("F23").Value = 6370
$ws.Range("F24").Value = 223
$ws.Range("F25").Value = 92
$ws.Range("F28").Value = 91
$ws.Range("F29").Value = 5761
$ws.Range("F30").Value = 1619
$ws.Range("F35").Value = 5900
$ws.Range("F36").Value = 101
$ws.Range("F39").Value = 77
$ws.Range("F40").Value = 4002
$ws.Range("F41").Value = 182
$ws.Range("F43").Value = 15
$ws.Range("F46").Value = 2394
$ws.Range("F49").Value = 1005
$ws.Range("F50").Value = 14
$ws.Range("F51").Value = 289

